# Rebuild the "Requisitos" bullet list: drop LOB1045 entirely, fix the
# "Àlgebra" -> "Álgebra" typo on LOB1037, and reorder the remaining
# course-requirement lines into their new sequence.

$d = $word.ActiveDocument

# Final desired order/content of the requirement lines.
$items = @(
  "LOQ4251 -  Fundamentos de Química  (Requisito)",
  "LOB1006 -  Cálculo IV  (Requisito)",
  "LOB1053 -  Física III  (Requisito)",
  "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito)",
  "LOB1003 -  Cálculo I  (Requisito)",
  "LOB1012 -  Estatística  (Requisito)",
  "LOB1024 -  Mecânica  (Requisito)",
  "LOB1036 -  Geometria Analítica  (Requisito)",
  "LOB1037 -  Álgebra Linear  (Requisito)",
  "LOB1038 -  Física Experimental I  (Requisito)",
  "LOB1039 -  Física Experimental III  (Requisito)",
  "LOB1041 -  Física Experimental II  (Requisito)",
  "LOB1052 -  Cálculo III  (Requisito)",
  "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)",
  "LOB1004 -  Cálculo II  (Requisito)",
  "LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito)",
  "LOB1018 -  Física I  (Requisito)",
  "LOB1019 -  Física II  (Requisito)"
)

# Locate the bullet-list paragraph holding the "(Requisito)" lines.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*(Requisito)*") {
        $target = $p
    }
}

$oldStart = $target.Range.Start
$oldEnd = $target.Range.End

# Manual line break character (same as a Word "Shift+Enter" <w:br/>).
$nl = [char]11

# Insert every new line, one at a time, right before the existing
# (soon to be removed) content. Each InsertBefore call on a point range
# produces its own run, so inserting in reverse order at the same fixed
# point reproduces the final forward order with one <w:r> per line.
for ($i = $items.Length - 1; $i -ge 0; $i--) {
    $insertPoint = $target.Range.Start
    $r = $d.Range($insertPoint, $insertPoint)
    $r.InsertBefore($items[$i] + $nl)
}

# Figure out how much text we just inserted, then remove everything
# after it up to the paragraph's original end -- that is the old set of
# requirement lines, now pushed later in the document.
$insertedLength = 0
foreach ($it in $items) {
    $insertedLength = $insertedLength + $it.Length + 1
}

$oldContent = $d.Range($oldStart + $insertedLength, $oldEnd + $insertedLength)
$oldContent.Text = ""
